$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.569.82"
$ws.Range("E2").Value = "  +7.49%  "
$ws.Range("D3").Value = "3.628.68"
$ws.Range("E3").Value = "  +7.50%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.89"
$ws.Range("E5").Value = "  +5.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.34"
$ws.Range("E6").Value = "  +9.59%  "
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("D8").Value = "3.604.89"
$ws.Range("E8").Value = "  +7.04%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.662"
$ws.Range("E11").Value = "  +4.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.63"
$ws.Range("E12").Value = "  +7.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000294"
$ws.Range("E13").Value = "  +6.30%  "
$ws.Range("E14").Value = "  +5.98%  "
$ws.Range("D15").Value = "4.208.26"
$ws.Range("E15").Value = "  +7.21%  "
$ws.Range("D16").Value = "3.630.05"
$ws.Range("E16").Value = "  +7.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.33"
$ws.Range("E17").Value = "  +6.38%  "
$ws.Range("D18").Value = "70.387.45"
$ws.Range("E18").Value = "  +7.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.58"
$ws.Range("E19").Value = "  +6.40%  "
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("E21").Value = "  +5.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.95"
$ws.Range("E22").Value = "  +5.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.54"
$ws.Range("E23").Value = "  +13.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.72"
$ws.Range("E24").Value = "  +17.50%  "
$ws.Range("E25").Value = "  +8.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.49"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("E27").Value = "  +6.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.21"
$ws.Range("E28").Value = "  +6.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.37"
$ws.Range("E29").Value = "  +7.79%  "
$ws.Range("E30").Value = "  +4.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.59"
$ws.Range("E31").Value = "  +15.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.23"
$ws.Range("E32").Value = "  +7.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "612.55"
$ws.Range("E33").Value = "  +6.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.29"
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("E35").Value = "  +7.68%  "
$ws.Range("D36").Value = "0.0₃0830"
$ws.Range("E36").Value = "  +12.53%  "
$ws.Range("E37").Value = "  +4.54%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.88"
$ws.Range("E39").Value = "  +5.80%  "
$ws.Range("E40").Value = "  +7.28%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").Value = "3.356.11"
$ws.Range("E42").Value = "  +8.47%  "
$ws.Range("E43").Value = "  +8.01%  "
$ws.Range("E44").Value = "  +7.19%  "
$ws.Range("E45").Value = "  +9.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  +6.82%  "
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.73"
$ws.Range("E48").Value = "  +11.03%  "
$ws.Range("E49").Value = "  +7.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.33"
$ws.Range("E50").Value = "  +6.32%  "
$ws.Range("E51").Value = "  -0.13%  "
